$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @('Belo Horizonte', 'Páprica Defumada Kitano Reserva 36g', 0, 'https://mercado.carrefour.com.br/farofa-de-mandioca-tradicional-yoki-400g-6582613/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/massa-para-pastel-discao-massa-leve-500g-841757/p', '2025-09-18'),
    @('Belo Horizonte', 'Pimentão Block Vermelho Trebeshi 150 g', 0, 'https://mercado.carrefour.com.br/pimentao-block-vermelho-trebeshi-150-g-5738458/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/chocolate-ao-leite-com-amendoim-shot-165g-5790859/p', '2025-09-18'),
    @('Belo Horizonte', 'Couve Flor Cledson 300 g', 0, 'https://mercado.carrefour.com.br/couve-flor-cledson-300-g-9560297/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/patinho-fracionado-a-vacuo-500g-18325/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/paleta-bovina-a-vacuo-500gnao-reativarcodigo-de-compra-20745/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/costela-minga-bovina-cong-aprox-2kg-224006/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/presunto-cozido-sem-capa-fatiado-aurora-aproximadamente-200-g-49450/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/mortadela-defumada-sadia-280g-5447045/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/queijo-minas-frescal-aurora-450-g-6264693/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/queijo-coalho-bom-leite-500-g-4305054/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/leite-uht-integral-piratininga-1-l-665017/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/iogurte-natural-tradicional-batavo-170g-5150439/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/manteiga-com-sal-aviacao-200-g-10010/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/creme-de-leite-ultrapasteurizado-itambe-200-g-5988921/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/requeijao-cremoso-aviacao-tradicional-220-g-10000/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/acucar-cristal-carrefour-1kg-5147300/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/mel-com-cacau-e-avela-400-g-4510146/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/geleia-de-goiaba-selecoes-c-pedacos-260-g-1280815/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/suco-de-uva-integral-maric-1-l-3538256/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/vinho-tinto-fino-seco-cabernet-sauvignon-pergola-750ml-1521709/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/whisky-red-label-johnnie-walker-1-litro-2719/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/refrigerante-coca-cola-sabor-cola-1-5-l-11087/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/cafe-torrado-e-moido-extraforte-melitta-500g-271203/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/farinha-de-trigo-dona-benta-tradicional-1kg-196416/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/azeite-extravirgem-portugues-oliveira-da-serra-500-ml-4526108/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/oleo-de-soja-soya-900ml-482616/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/margarina-qualy-com-sal-250g-4815618/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/arroz-branco-longofino-tipo-1-tio-joao-1kg-115658/p', '2025-09-18'),
    @('Belo Horizonte', 'Não encontrado', 0, 'https://mercado.carrefour.com.br/feijao-preto-tipo-1-kicaldo-1kg-466510/p', '2025-09-18'),
)

$startRow = 33
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[4]
    $eCell.ClearFormats()
}
